$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimina EC anteriores y se agregan nuevos" - the previous "Periodo Mora"
# values (1812, 1901, 1902, 1903) in column E are replaced with the new
# period set, written in this order so new shared-string entries are
# appended as 1903, 1902, 1901 (in that order) while 1812 is reused.
$ws.Range("E16").Value = "1903"
$ws.Range("E17").Value = "1902"
$ws.Range("E18").Value = "1901"
$ws.Range("E19").Value = "1812"
